$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(49, "face", "顔|かお")
    ,@(50, "emoticon", "顔文字|かおもじ")
    ,@(51, "complexion", "顔色|かおいろ")
    ,@(52, "smiling face", "笑顔|えがお")
    ,@(53, "washing one's face", "洗顔|せんがん")
    ,@(54, "sad", "悲しい|かなしい")
    ,@(55, "tragedy", "悲劇|ひげき")
    ,@(56, "miserable", "悲惨な|ひさんな")
    ,@(57, "to grieve", "悲しむ|かなしむ")
    ,@(58, "to get angry", "怒る|おこる")
    ,@(59, "anger; rage", "怒り|いかり")
    ,@(60, "human emotions", "喜怒哀楽|きどあいらく")
    ,@(61, "different", "違う|ちがう")
    ,@(62, "difference", "違い|ちがい")
    ,@(63, "to make a mistake", "間違える|まちがえる")
    ,@(64, "illegal", "違法|いほう")
    ,@(65, "violation", "違反|いはん")
    ,@(66, "strange", "変な|へんな")
    ,@(67, "tough; hectic", "大変な|たいへんな")
    ,@(68, "change", "変化|へんか")
    ,@(69, "eccentric person", "変人|へんじん")
    ,@(70, "to change (something)", "変える|かえる")
    ,@(71, "to compare", "比べる|くらべる")
    ,@(72, "comparison", "比較|ひかく")
    ,@(73, "proportion", "比例|ひれい")
    ,@(74, "contrast", "対比|たいひ")
    ,@(75, "figure of speech", "比喩|ひゆ")
    ,@(76, "expression", "表情|ひょうじょう")
    ,@(77, "to sympathize", "同情する|どうじょうする")
    ,@(78, "friendship", "友情|ゆうじょう")
    ,@(79, "information", "情報|じょうほう")
    ,@(80, "mercy", "情け|なさけ")
    ,@(81, "emotion", "感情|かんじょう")
    ,@(82, "to be moved", "感動する|かんどうする")
    ,@(83, "to feel", "感じる|かんじる")
    ,@(84, "gratitude", "感謝|かんしゃ")
    ,@(85, "impression", "感想|かんそう")
    ,@(86, "survey; research", "調査|ちょうさ")
    ,@(87, "to look into; to examine", "調べる|しらべる")
    ,@(88, "condition", "調子|ちょうし")
    ,@(89, "to emphasize", "強調する|きょうちょうする")
    ,@(90, "survey; research", "調査|ちょうさ")
    ,@(91, "inspection", "検査|けんさ")
    ,@(92, "screening", "審査|しんさ")
    ,@(93, "criminal investigation", "捜査|そうさ")
    ,@(94, "result", "結果|けっか")
    ,@(95, "fruit", "果物|くだもの")
    ,@(96, "fruit juice", "果汁|かじゅう")
    ,@(97, "effect", "効果|こうか")
    ,@(98, "to use up", "使い果たす|つかいはたす")
    ,@(99, "culture", "文化|ぶんか")
    ,@(100, "chemistry", "化学|かがく")
    ,@(101, "assimilation", "同化|どうか")
    ,@(102, "goblin; ghost", "お化け|おばけ")
    ,@(103, "makeup", "化粧|けしょう")
    ,@(104, "side", "横|よこ")
    ,@(105, "horizontal writing", "横書き|よこがき")
    ,@(106, "grand champion of sumo", "横綱|よこづな")
    ,@(107, "to traverse", "横断する|おうだんする")
    ,@(108, "partner", "相手|あいて")
    ,@(109, "prime minister", "首相|しゅしょう")
    ,@(110, "consultation", "相談|そうだん")
    ,@(111, "mutual", "相互の|そうごの")
    ,@(112, "to answer", "答える|こたえる")
    ,@(113, "answer", "答え／答|こたえ")
    ,@(114, "reply; answer", "回答|かいとう")
    ,@(115, "correct answer", "正答|せいとう")
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
}
